$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting the existing rows 44:184 down to 45:185.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the same "reference" data that the
# (now shifted-down) old row 44 carries, except for the Fecha (D) and
# Volumen (J) columns which get new values.
$ws.Cells.Item(44, 1).Value = 5
$ws.Cells.Item(44, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value = "Maule"
$ws.Cells.Item(44, 4).Value = 44481
$ws.Cells.Item(44, 5).Value = 7
$ws.Cells.Item(44, 6).Value = 100114014
$ws.Cells.Item(44, 7).Value = "Betarraga"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 4000
$ws.Cells.Item(44, 11).Value = 600
$ws.Cells.Item(44, 12).Value = 600
$ws.Cells.Item(44, 13).Value = 600
$ws.Cells.Item(44, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 120
$ws.Cells.Item(44, 17).Value = 5
$ws.Cells.Item(44, 18).Value = "Hortaliza"
